# Generate Report for Handback
# Adds a new handback row (b512c95f-9fa1-403b-a53b-5c03f44ede5c) to the
# Overview / zh-cn / de-de sheets, mirroring the existing
# "751d3533-e448-4877-9d0d-898286d1e54e" ("in sync with en-US") entries.

$wb = $excel.ActiveWorkbook

$fileGuid  = "b512c95f-9fa1-403b-a53b-5c03f44ede5c"
$fileName  = "$fileGuid.md"
$display   = "e2e\$fileGuid.md"
$status    = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$rowOv = $wsOverview.ListObjects.Item(1).ListRows.Add().Range.Row

$wsOverview.Cells.Item($rowOv, 1).Value = $fileName
$wsOverview.Cells.Item($rowOv, 2).Value = $display
$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item($rowOv, 2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c31a176338e3cf3b6dad538daf208ee736c06c10/e2e/$fileName", "", "", $display)
$wsOverview.Cells.Item($rowOv, 3).Value = ".md"
$wsOverview.Cells.Item($rowOv, 5).Value = $status
$wsOverview.Cells.Item($rowOv, 6).Value = $status
$wsOverview.Cells.Item($rowOv, 7).Value = "2016-09-02 22:49:45"
$wsOverview.Cells.Item($rowOv, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$rowZh = $wsZh.ListObjects.Item(1).ListRows.Add().Range.Row

$xlfZh = "$fileGuid.b7cc76f1e5665da7fd6316deca7b8037124fce0d.zh-cn.xlf"

$wsZh.Cells.Item($rowZh, 1).Value = $fileName
$wsZh.Hyperlinks.Add($wsZh.Cells.Item($rowZh, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c31a176338e3cf3b6dad538daf208ee736c06c10/e2e/$fileName", "", "", $fileName)
$wsZh.Cells.Item($rowZh, 2).Value = ".md"
$wsZh.Cells.Item($rowZh, 3).Value = $status
$wsZh.Cells.Item($rowZh, 4).Value = "e2e"
$wsZh.Cells.Item($rowZh, 5).Value = "ht"
$wsZh.Cells.Item($rowZh, 6).Value = "True"
$wsZh.Cells.Item($rowZh, 7).Value = $xlfZh
$wsZh.Cells.Item($rowZh, 8).Value = "2016-09-02 22:49:41"
$wsZh.Cells.Item($rowZh, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item($rowZh, 9).Value = $fileName
$wsZh.Hyperlinks.Add($wsZh.Cells.Item($rowZh, 9), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/74d29fc30fa453aa605745d4b6ab154942b3a8a6/e2e/$fileName", "", "", $fileName)
$wsZh.Cells.Item($rowZh, 10).Value = $xlfZh
$wsZh.Cells.Item($rowZh, 11).Value = "2016-09-02 22:50:07"
$wsZh.Cells.Item($rowZh, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item($rowZh, 13).Value = "True"
$wsZh.Cells.Item($rowZh, 15).Value = "False"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$rowDe = $wsDe.ListObjects.Item(1).ListRows.Add().Range.Row

$xlfDe = "$fileGuid.b7cc76f1e5665da7fd6316deca7b8037124fce0d.de-de.xlf"

$wsDe.Cells.Item($rowDe, 1).Value = $fileName
$wsDe.Hyperlinks.Add($wsDe.Cells.Item($rowDe, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c31a176338e3cf3b6dad538daf208ee736c06c10/e2e/$fileName", "", "", $fileName)
$wsDe.Cells.Item($rowDe, 2).Value = ".md"
$wsDe.Cells.Item($rowDe, 3).Value = $status
$wsDe.Cells.Item($rowDe, 4).Value = "e2e"
$wsDe.Cells.Item($rowDe, 5).Value = "ht"
$wsDe.Cells.Item($rowDe, 6).Value = "True"
$wsDe.Cells.Item($rowDe, 7).Value = $xlfDe
$wsDe.Cells.Item($rowDe, 8).Value = "2016-09-02 22:49:45"
$wsDe.Cells.Item($rowDe, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item($rowDe, 9).Value = $fileName
$wsDe.Hyperlinks.Add($wsDe.Cells.Item($rowDe, 9), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/682ae117fd55855bcf30f74e334468718c7d0927/e2e/$fileName", "", "", $fileName)
$wsDe.Cells.Item($rowDe, 10).Value = $xlfDe
$wsDe.Cells.Item($rowDe, 11).Value = "2016-09-02 22:50:18"
$wsDe.Cells.Item($rowDe, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item($rowDe, 13).Value = "True"
$wsDe.Cells.Item($rowDe, 15).Value = "False"
